# Refresh the cryptos list (Price / Volume(1h) columns) with latest scraped
# values, as produced by the scheduled GitHub Actions scrape job.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Note: some new Price values (col D) look like plain decimals (e.g. "1.00",
# "601.01") which Excel would otherwise auto-convert to numbers, dropping the
# trailing zero / formatting that the source keeps as text. Writing them with
# a leading apostrophe forces text, then resetting the style back to "Normal"
# clears the quote-prefix formatting flag so the cell style stays untouched
# (matching the rest of the sheet, which never used the quote-prefix style).

$ws.Range("D2").Value = "68.886.63"
$ws.Range("E2").Value = "  -0.37%  "
$ws.Range("D3").Value = "3.862.86"
$ws.Range("E3").Value = "  +2.82%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "'601.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'162.48"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.80%  "
$ws.Range("D7").Value = "3.862.69"
$ws.Range("E7").Value = "  +2.85%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -1.55%  "
$ws.Range("E10").Value = "  -0.90%  "
$ws.Range("D11").Value = "'6.31"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.50%  "
$ws.Range("E12").Value = "  -0.15%  "
$ws.Range("E13").Value = "  -2.90%  "
$ws.Range("E14").Value = "  -2.06%  "
$ws.Range("D15").Value = "4.510.52"
$ws.Range("E15").Value = "  +2.85%  "
$ws.Range("D16").Value = "3.845.37"
$ws.Range("E16").Value = "  +2.49%  "
$ws.Range("D17").Value = "69.056.13"
$ws.Range("E17").Value = "  -0.16%  "
$ws.Range("E18").Value = "  +2.51%  "
$ws.Range("E19").Value = "  -0.31%  "
$ws.Range("D20").Value = "'11.37"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.25%  "
$ws.Range("E21").Value = "  -1.33%  "
$ws.Range("D22").Value = "'485.03"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.73%  "
$ws.Range("E23").Value = "  -1.49%  "
$ws.Range("E24").Value = "  +7.09%  "
$ws.Range("D25").Value = "'84.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.13%  "
$ws.Range("E26").Value = "  -2.87%  "
$ws.Range("D27").Value = "'12.11"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.63%  "
$ws.Range("E28").Value = "  -0.07%  "
$ws.Range("D29").Value = "'9.97"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.40%  "
$ws.Range("E30").Value = "  -1.06%  "
$ws.Range("E31").Value = "  -3.47%  "
$ws.Range("D32").Value = "4.011.01"
$ws.Range("E32").Value = "  +2.76%  "
$ws.Range("D33").Value = "'32.33"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.32%  "
$ws.Range("E34").Value = "  -4.35%  "
$ws.Range("D35").Value = "3.809.67"
$ws.Range("E35").Value = "  +3.26%  "
$ws.Range("E36").Value = "  -1.51%  "
$ws.Range("E37").Value = "  +1.09%  "
$ws.Range("E38").Value = "  +2.00%  "
$ws.Range("D39").Value = "'5.89"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.59%  "
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("E41").Value = "  -2.38%  "
$ws.Range("D42").Value = "'437.28"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.71%  "
$ws.Range("E43").Value = "  -2.18%  "
$ws.Range("E44").Value = "  -0.30%  "
$ws.Range("E45").Value = "  -0.02%  "
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("E47").Value = "  -0.82%  "
$ws.Range("D48").Value = "'143.44"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.57%  "
$ws.Range("D49").Value = "2.839.32"
$ws.Range("E49").Value = "  +1.55%  "
$ws.Range("D50").Value = "'0.0358"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.40%  "
$ws.Range("D51").Value = "'26.06"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +13.16%  "
